$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @('Dollar', '4,91', '12:29', ' sábado, 15 de abril de 2023 ')
$data[3] = @('Euro', '5,45', '12:29', ' sábado, 15 de abril de 2023 ')
$data[4] = @('Iene', '0,037', '12:29', ' sábado, 15 de abril de 2023 ')
$data[5] = @('Yuan Chinês', '0,71', '12:29', ' sábado, 15 de abril de 2023 ')
$data[6] = @('Dollar', '4,93', '18:31', ' quinta-feira, 13 de abril de 2023 ')
$data[7] = @('Euro', '5,44', '18:31', ' quinta-feira, 13 de abril de 2023 ')
$data[8] = @('Iene', '0,037', '18:31', ' quinta-feira, 13 de abril de 2023 ')
$data[9] = @('Yuan Chinês', '0,72', '18:31', ' quinta-feira, 13 de abril de 2023 ')
$data[10] = @('Dollar', '4,93', '17:50', ' quinta-feira, 13 de abril de 2023 ')
$data[11] = @('Euro', '5,45', '17:50', ' quinta-feira, 13 de abril de 2023 ')
$data[12] = @('Iene', '0,037', '17:50', ' quinta-feira, 13 de abril de 2023 ')
$data[13] = @('Yuan Chinês', '0,72', '17:50', ' quinta-feira, 13 de abril de 2023 ')
$data[14] = @('Dollar', '4,93', '17:48', ' quinta-feira, 13 de abril de 2023 ')
$data[15] = @('Euro', '5,45', '17:48', ' quinta-feira, 13 de abril de 2023 ')
$data[16] = @('Iene', '0,037', '17:48', ' quinta-feira, 13 de abril de 2023 ')
$data[17] = @('Yuan Chinês', '0,72', '17:48', ' quinta-feira, 13 de abril de 2023 ')
$data[18] = @('Dollar', '4,93', '17:47', ' quinta-feira, 13 de abril de 2023 ')
$data[19] = @('Euro', '5,44', '17:47', ' quinta-feira, 13 de abril de 2023 ')
$data[20] = @('Iene', '0,037', '17:47', ' quinta-feira, 13 de abril de 2023 ')
$data[21] = @('Yuan Chinês', '0,72', '17:47', ' quinta-feira, 13 de abril de 2023 ')
$data[22] = @('Dollar', '4,93', '17:45', ' quinta-feira, 13 de abril de 2023 ')
$data[23] = @('Euro', '5,44', '17:45', ' quinta-feira, 13 de abril de 2023 ')
$data[24] = @('Iene', '0,037', '17:45', ' quinta-feira, 13 de abril de 2023 ')
$data[25] = @('Yuan Chinês', '0,72', '17:45', ' quinta-feira, 13 de abril de 2023 ')
$data[26] = @('Dollar', '4,93', '17:38', ' quinta-feira, 13 de abril de 2023 ')
$data[27] = @('Euro', '5,45', '17:38', ' quinta-feira, 13 de abril de 2023 ')
$data[28] = @('Iene', '0,037', '17:38', ' quinta-feira, 13 de abril de 2023 ')
$data[29] = @('Yuan Chinês', '0,72', '17:38', ' quinta-feira, 13 de abril de 2023 ')
$data[30] = @('Dollar', '4,93', '17:32', ' quinta-feira, 13 de abril de 2023 ')
$data[31] = @('Euro', '5,44', '17:32', ' quinta-feira, 13 de abril de 2023 ')
$data[32] = @('Iene', '0,037', '17:32', ' quinta-feira, 13 de abril de 2023 ')
$data[33] = @('Yuan Chinês', '0,72', '17:32', ' quinta-feira, 13 de abril de 2023 ')
$data[34] = @('Dollar', '4,93', '17:29', ' quinta-feira, 13 de abril de 2023 ')
$data[35] = @('Euro', '5,44', '17:29', ' quinta-feira, 13 de abril de 2023 ')
$data[36] = @('Iene', '0,037', '17:29', ' quinta-feira, 13 de abril de 2023 ')
$data[37] = @('Yuan Chinês', '0,72', '17:29', ' quinta-feira, 13 de abril de 2023 ')
$data[38] = @('Dollar', '5,03', '22:32', ' quarta-feira, 5 de abril de 2023 ')
$data[39] = @('Euro', '5,48', '22:32', ' quarta-feira, 5 de abril de 2023 ')
$data[40] = @('Iene', '0,038', '22:32', ' quarta-feira, 5 de abril de 2023 ')
$data[41] = @('Yuan Chinês', '0,73', '22:32', ' quarta-feira, 5 de abril de 2023 ')
$data[42] = @('Dollar', '5,03', '21:21', ' quarta-feira, 5 de abril de 2023 ')
$data[43] = @('Euro', '5,49', '21:21', ' quarta-feira, 5 de abril de 2023 ')
$data[44] = @('Iene', '0,038', '21:21', ' quarta-feira, 5 de abril de 2023 ')
$data[45] = @('Yuan Chinês', '0,73', '21:21', ' quarta-feira, 5 de abril de 2023 ')
$data[46] = @('Dollar', '5,03', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[47] = @('Euro', '5,49', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[48] = @('Iene', '0,038', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[49] = @('Yuan Chinês', '0,73', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[50] = @('Dollar', '5,03', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[51] = @('Euro', '5,49', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[52] = @('Iene', '0,038', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[53] = @('Yuan Chinês', '0,73', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[54] = @('Dollar', '4,93', '17:38', ' quinta-feira, 13 de abril de 2023 ')
$data[55] = @('Euro', '5,45', '17:38', ' quinta-feira, 13 de abril de 2023 ')
$data[56] = @('Iene', '0,037', '17:38', ' quinta-feira, 13 de abril de 2023 ')
$data[57] = @('Yuan Chinês', '0,72', '17:38', ' quinta-feira, 13 de abril de 2023 ')
$data[58] = @('Dollar', '4,93', '17:32', ' quinta-feira, 13 de abril de 2023 ')
$data[59] = @('Euro', '5,44', '17:32', ' quinta-feira, 13 de abril de 2023 ')
$data[60] = @('Iene', '0,037', '17:32', ' quinta-feira, 13 de abril de 2023 ')
$data[61] = @('Yuan Chinês', '0,72', '17:32', ' quinta-feira, 13 de abril de 2023 ')
$data[62] = @('Dollar', '4,93', '17:29', ' quinta-feira, 13 de abril de 2023 ')
$data[63] = @('Euro', '5,44', '17:29', ' quinta-feira, 13 de abril de 2023 ')
$data[64] = @('Iene', '0,037', '17:29', ' quinta-feira, 13 de abril de 2023 ')
$data[65] = @('Yuan Chinês', '0,72', '17:29', ' quinta-feira, 13 de abril de 2023 ')
$data[66] = @('Dollar', '5,03', '22:32', ' quarta-feira, 5 de abril de 2023 ')
$data[67] = @('Euro', '5,48', '22:32', ' quarta-feira, 5 de abril de 2023 ')
$data[68] = @('Iene', '0,038', '22:32', ' quarta-feira, 5 de abril de 2023 ')
$data[69] = @('Yuan Chinês', '0,73', '22:32', ' quarta-feira, 5 de abril de 2023 ')
$data[70] = @('Dollar', '5,03', '21:21', ' quarta-feira, 5 de abril de 2023 ')
$data[71] = @('Euro', '5,49', '21:21', ' quarta-feira, 5 de abril de 2023 ')
$data[72] = @('Iene', '0,038', '21:21', ' quarta-feira, 5 de abril de 2023 ')
$data[73] = @('Yuan Chinês', '0,73', '21:21', ' quarta-feira, 5 de abril de 2023 ')
$data[74] = @('Dollar', '5,03', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[75] = @('Euro', '5,49', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[76] = @('Iene', '0,038', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[77] = @('Yuan Chinês', '0,73', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[78] = @('Dollar', '5,03', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[79] = @('Euro', '5,49', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[80] = @('Iene', '0,038', '21:20', ' quarta-feira, 5 de abril de 2023 ')
$data[81] = @('Yuan Chinês', '0,73', '21:20', ' quarta-feira, 5 de abril de 2023 ')

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item($r, 1).Value = ($r - 2)
    $ws.Cells.Item($r, 2).Value = $row[0]
    $cCell = $ws.Cells.Item($r, 3)
    if ($row[1] -eq '0,037' -or $row[1] -eq '0,038') {
        $cCell.NumberFormat = "@"
    }
    $cCell.Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A78:A81").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Output "applied edits"
